$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.400.92"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.934.41"
$ws.Range("E3").Value = "  -2.69%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.82"
$ws.Range("E5").Value = "  -3.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.91"
$ws.Range("E6").Value = "  +1.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.931.35"
$ws.Range("E9").Value = "  -2.65%  "

$ws.Range("E10").Value = "  -3.95%  "

$ws.Range("E11").Value = "  -3.72%  "

$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("E13").Value = "  -1.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.23"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.362.65"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.422.71"
$ws.Range("E17").Value = "  -2.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.90"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.939.44"
$ws.Range("E19").Value = "  -2.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.01"
$ws.Range("E20").Value = "  +8.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "444.43"
$ws.Range("E21").Value = "  -4.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.685"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.24"
$ws.Range("E23").Value = "  -1.76%  "

$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.11"
$ws.Range("E26").Value = "  -3.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("E27").Value = "  -6.60%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").Value = "  +1.95%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("E31").Value = "  -1.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000101"
$ws.Range("E32").Value = "  -5.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.05"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("E34").Value = "  -1.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.970"
$ws.Range("E36").Value = "  -2.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.70"
$ws.Range("E37").Value = "  -1.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.68"
$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "44.81"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("E40").Value = "  -9.84%  "

$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "384.79"
$ws.Range("E45").Value = "  -2.37%  "

$ws.Range("E46").Value = "  -1.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.699.83"
$ws.Range("E47").Value = "  -3.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.25"
$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("E50").Value = "  +4.16%  "

# Row 41: dogwifhat -> Kaspa
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.120"
$ws.Range("E41").Value = "  -2.66%  "

# Row 43: Kaspa -> dogwifhat
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.83"
$ws.Range("E43").Value = "  -7.68%  "

# Row 51: Stellar -> InjectiveProtocol
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.33"
$ws.Range("E51").Value = "  -1.48%  "
